# Implemented Status tracking column in TransactionData to conform to solution standards
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Constants")

$ws.Range("A12").Value = "Status_Success"
$ws.Range("B12").Value = "Success"
$ws.Range("C12").Value = "Status message for TransactionData to record successful transaction."

$ws.Range("A13").Value = "Status_Failure"
$ws.Range("B13").Value = "Failed"
$ws.Range("C13").Value = "Status message for TransactionData to record failed transaction."

$ws.Range("A14").Value = "Status_Pending"
$ws.Range("B14").Value = "Pending"
$ws.Range("C14").Value = "Status message for TransactionData to record pending transaction."

$ws.Activate()
$ws.Range("A12:C14").Select()
